$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free direct updates. Columns B/C/E are plain text (names, URLs,
# percentage strings) so a normal .Value assignment keeps them as text.
# Column D holds numeric-looking text (prices) that Excel would otherwise
# auto-convert to a real number (dropping formatting like trailing zeros),
# so those cells are written via a temporary Text number-format, then
# ClearFormats() removes the temporary formatting again (style stays default).

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '37.379.80'
$c.ClearFormats()

$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.098.32'
$c.ClearFormats()
$ws.Range('E3').Value = '  +3.78%  '

$ws.Range('E4').Value = '  +0.03%  '

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '250.31'
$c.ClearFormats()
$ws.Range('E5').Value = '  +1.21%  '

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.660'
$c.ClearFormats()
$ws.Range('E6').Value = '  -0.30%  '

$ws.Range('E7').Value = '  +0.02%  '

$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '51.92'
$c.ClearFormats()
$ws.Range('E8').Value = '  +13.90%  '

$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '62.03'
$c.ClearFormats()
$ws.Range('E9').Value = '  +8.86%  '

$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.374'
$c.ClearFormats()
$ws.Range('E10').Value = '  +2.71%  '

$ws.Range('E11').Value = '  +3.32%  '

$ws.Range('E12').Value = '  +6.76%  '

$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '15.09'
$c.ClearFormats()
$ws.Range('E13').Value = '  +2.39%  '

$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '2.400.77'
$c.ClearFormats()
$ws.Range('E14').Value = '  +3.68%  '

$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '0.834'
$c.ClearFormats()
$ws.Range('E15').Value = '  +3.12%  '

$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '2.092.32'
$c.ClearFormats()
$ws.Range('E16').Value = '  +3.63%  '

$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '5.12'
$c.ClearFormats()
$ws.Range('E17').Value = '  +4.05%  '

$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '37.398.42'
$c.ClearFormats()
$ws.Range('E18').Value = '  +2.55%  '

$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '72.16'
$c.ClearFormats()
$ws.Range('E19').Value = '  +1.28%  '

$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '13.96'
$c.ClearFormats()
$ws.Range('E20').Value = '  +7.20%  '

$ws.Range('E21').Value = '  +1.39%  '

$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '240.45'
$c.ClearFormats()
$ws.Range('E22').Value = '  +2.40%  '

$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '5.23'
$c.ClearFormats()
$ws.Range('E23').Value = '  +4.86%  '

$ws.Range('E24').Value = '  -0.07%  '

$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '2.45'
$c.ClearFormats()
$ws.Range('E25').Value = '  -1.12%  '

$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '171.17'
$c.ClearFormats()
$ws.Range('E26').Value = '  +4.69%  '

$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '9.20'
$c.ClearFormats()
$ws.Range('E27').Value = '  +7.51%  '

$ws.Range('E28').Value = '  +5.72%  '

$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '2.00'
$c.ClearFormats()
$ws.Range('E29').Value = '  -0.46%  '

$ws.Range('E30').Value = '  +0.18%  '

$ws.Range('E31').Value = '  +25.31%  '

$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '4.49'
$c.ClearFormats()
$ws.Range('E32').Value = '  +1.93%  '

$ws.Range('E33').Value = '  +3.38%  '

$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '0.0914'
$c.ClearFormats()
$ws.Range('E34').Value = '  +8.87%  '

$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '20.27'
$c.ClearFormats()
$ws.Range('E35').Value = '  -4.77%  '

$ws.Range('E36').Value = '  -0.08%  '

$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '2.32'
$c.ClearFormats()
$ws.Range('E37').Value = '  +6.53%  '

$ws.Range('E38').Value = '  -0.61%  '

$ws.Range('E39').Value = '  +1.16%  '

$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '1.33'
$c.ClearFormats()
$ws.Range('E40').Value = '  -2.56%  '

$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '18.23'
$c.ClearFormats()
$ws.Range('E41').Value = '  +12.32%  '

$ws.Range('E42').Value = '  +3.66%  '

$ws.Range('E43').Value = '  +7.19%  '

$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '99.74'
$c.ClearFormats()
$ws.Range('E44').Value = '  +2.90%  '

$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '0.0908'
$c.ClearFormats()
$ws.Range('E45').Value = '  +11.88%  '

$ws.Range('B46').Value = 'HuobiToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '2.75'
$c.ClearFormats()
$ws.Range('E46').Value = '  +0.33%  '

$ws.Range('E47').Value = '  +8.15%  '

$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '1.321.16'
$c.ClearFormats()
$ws.Range('E48').Value = '  +0.33%  '

$ws.Range('E49').Value = '  +14.29%  '

$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '2.296.37'
$c.ClearFormats()
$ws.Range('E50').Value = '  +3.41%  '

$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '2.27'
$c.ClearFormats()
$ws.Range('E51').Value = '  +1.45%  '
